$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 233.54546
$v = $ws.Range("H2").Value2
Write-Host $v
